# Auto-generated: scheduled market-data refresh for Omega_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 134086
$ws.Range("J28").Value = 564.375
$ws.Range("L28").Value = 564.375
$ws.Range("N28").Value = -1534.375
$ws.Range("H98").Value = 1075.5264
$ws.Range("I98").Value = 1066.0667
$ws.Range("J98").Value = 1111
$ws.Range("K98").Value = 1066.0667
$ws.Range("L98").Value = 1111
$ws.Range("M98").Value = 431.9332999999999
$ws.Range("N98").Value = -4107
$ws.Range("H112").Value = 3285
$ws.Range("J112").Value = 3352.9
$ws.Range("L112").Value = 10058.7
$ws.Range("N112").Value = -12274.7
$ws.Range("H122").Value = 1075.5264
$ws.Range("I122").Value = 1066.0667
$ws.Range("J122").Value = 1111
$ws.Range("K122").Value = 3198.2001
$ws.Range("L122").Value = 3333
$ws.Range("M122").Value = -748.2001
$ws.Range("N122").Value = -8233
$ws.Range("H137").Value = 2492.889
$ws.Range("I137").Value = 1750
$ws.Range("J137").Value = 2864.3333
$ws.Range("K137").Value = 5250
$ws.Range("L137").Value = 8592.999899999999
$ws.Range("M137").Value = -2700
$ws.Range("N137").Value = -13692.9999
$ws.Range("H138").Value = 8442.212
$ws.Range("I138").Value = 4219.5713
$ws.Range("J138").Value = 9579.076999999999
$ws.Range("K138").Value = 12658.7139
$ws.Range("L138").Value = 28737.231
$ws.Range("M138").Value = -7518.713899999999
$ws.Range("N138").Value = -39017.231

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1117.8334
$ws.Range("I5").Value = 1117.8334
$ws.Range("K5").Value = 1117.8334
$ws.Range("M5").Value = -1005.8334
$ws.Range("H32").Value = 6311.5713
$ws.Range("I32").Value = 645.98114
$ws.Range("J32").Value = 36339.2
$ws.Range("K32").Value = 645.98114
$ws.Range("L32").Value = 36339.2
$ws.Range("M32").Value = -358.98114
$ws.Range("N32").Value = -36913.2
$ws.Range("H45").Value = 3385.3215
$ws.Range("I45").Value = 2698.6875
$ws.Range("K45").Value = 2698.6875
$ws.Range("M45").Value = -2321.6875
$ws.Range("H61").Value = 5305.4346
$ws.Range("I61").Value = 5005.952
$ws.Range("J61").Value = 8450
$ws.Range("K61").Value = 5005.952
$ws.Range("L61").Value = 8450
$ws.Range("M61").Value = -4793.952
$ws.Range("N61").Value = -8874
$ws.Range("H74").Value = 2116.7646
$ws.Range("I74").Value = 1585.7142
$ws.Range("K74").Value = 1585.7142
$ws.Range("M74").Value = -711.7141999999999
$ws.Range("H77").Value = 2116.7646
$ws.Range("I77").Value = 1585.7142
$ws.Range("K77").Value = 7928.571
$ws.Range("M77").Value = -3560.571
$ws.Range("H94").Value = 34333
$ws.Range("J94").Value = 34333
$ws.Range("L94").Value = 34333
$ws.Range("N94").Value = -36135
$ws.Range("H122").Value = 1888.5883
$ws.Range("I122").Value = 1842.0714
$ws.Range("J122").Value = 2105.6667
$ws.Range("K122").Value = 5526.2142
$ws.Range("L122").Value = 6317.000100000001
$ws.Range("M122").Value = -3076.2142
$ws.Range("N122").Value = -11217.0001
$ws.Range("H136").Value = 5305.4346
$ws.Range("I136").Value = 5005.952
$ws.Range("J136").Value = 8450
$ws.Range("K136").Value = 15017.856
$ws.Range("L136").Value = 25350
$ws.Range("M136").Value = -12467.856
$ws.Range("N136").Value = -30450

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1117.8334
$ws.Range("I4").Value = 1117.8334
$ws.Range("K4").Value = 1117.8334
$ws.Range("M4").Value = -1002.8334
$ws.Range("H86").Value = 1776.2916
$ws.Range("I86").Value = 1381.5883
$ws.Range("J86").Value = 2734.8572
$ws.Range("K86").Value = 1381.5883
$ws.Range("L86").Value = 2734.8572
$ws.Range("M86").Value = -258.5882999999999
$ws.Range("N86").Value = -4980.8572
$ws.Range("H89").Value = 1776.2916
$ws.Range("I89").Value = 1381.5883
$ws.Range("J89").Value = 2734.8572
$ws.Range("K89").Value = 6907.941499999999
$ws.Range("L89").Value = 13674.286
$ws.Range("M89").Value = -1291.941499999999
$ws.Range("N89").Value = -24906.286
$ws.Range("H105").Value = 3414.6428
$ws.Range("I105").Value = 2523.4614
$ws.Range("J105").Value = 15000
$ws.Range("K105").Value = 2523.4614
$ws.Range("L105").Value = 15000
$ws.Range("M105").Value = -776.4614000000001
$ws.Range("N105").Value = -18494

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 543608.5600000001
$ws.Range("J4").Value = 300876.66
$ws.Range("L4").Value = 300876.66
$ws.Range("N4").Value = -301100.66
$ws.Range("H31").Value = 4848.447
$ws.Range("I31").Value = 6910.6
$ws.Range("J31").Value = 3881.8125
$ws.Range("K31").Value = 6910.6
$ws.Range("L31").Value = 3881.8125
$ws.Range("M31").Value = -6615.6
$ws.Range("N31").Value = -4471.8125
$ws.Range("H34").Value = 4848.447
$ws.Range("I34").Value = 6910.6
$ws.Range("J34").Value = 3881.8125
$ws.Range("K34").Value = 6910.6
$ws.Range("L34").Value = 3881.8125
$ws.Range("M34").Value = -6708.6
$ws.Range("N34").Value = -4285.8125
$ws.Range("H105").Value = 2541.5
$ws.Range("I105").Value = 2388.6667
$ws.Range("K105").Value = 2388.6667
$ws.Range("M105").Value = -641.6667000000002
$ws.Range("H132").Value = 6252.7144
$ws.Range("I132").Value = 5602.1816
$ws.Range("J132").Value = 8638
$ws.Range("K132").Value = 16806.5448
$ws.Range("L132").Value = 25914
$ws.Range("M132").Value = -14276.5448
$ws.Range("N132").Value = -30974
$ws.Range("H134").Value = 7403.636
$ws.Range("I134").Value = 5740.6665
$ws.Range("K134").Value = 17221.9995
$ws.Range("M134").Value = -14686.9995
$ws.Range("H138").Value = 82139.42999999999
$ws.Range("J138").Value = 82139.42999999999
$ws.Range("L138").Value = 82139.42999999999
$ws.Range("N138").Value = -92419.42999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 47625844
$ws.Range("I8").Value = 47625844
$ws.Range("K8").Value = 142877532
$ws.Range("M8").Value = -142877393
$ws.Range("H23").Value = 1349.8334
$ws.Range("I23").Value = 899.75
$ws.Range("J23").Value = 2250
$ws.Range("K23").Value = 2699.25
$ws.Range("L23").Value = 6750
$ws.Range("M23").Value = -2464.25
$ws.Range("N23").Value = -7220
$ws.Range("H38").Value = 93.8
$ws.Range("I38").Value = 83.14286
$ws.Range("J38").Value = 107.36364
$ws.Range("K38").Value = 249.42858
$ws.Range("L38").Value = 322.09092
$ws.Range("M38").Value = 97.57141999999999
$ws.Range("N38").Value = -1016.09092
$ws.Range("H68").Value = 2894.0557
$ws.Range("I68").Value = 2660
$ws.Range("J68").Value = 2923.3125
$ws.Range("K68").Value = 7980
$ws.Range("L68").Value = 8769.9375
$ws.Range("M68").Value = -7169
$ws.Range("N68").Value = -10391.9375
$ws.Range("H71").Value = 2894.0557
$ws.Range("I71").Value = 2660
$ws.Range("J71").Value = 2923.3125
$ws.Range("K71").Value = 23940
$ws.Range("L71").Value = 26309.8125
$ws.Range("M71").Value = -19884
$ws.Range("N71").Value = -34421.8125
$ws.Range("H93").Value = 6641.6665
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 8962.5
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 26887.5
$ws.Range("M93").Value = -4128
$ws.Range("N93").Value = -30631.5
$ws.Range("H128").Value = 159504.75
$ws.Range("I128").Value = 159504.75
$ws.Range("K128").Value = 478514.25
$ws.Range("M128").Value = -473534.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2073.578
$ws.Range("I122").Value = 1205.9429
$ws.Range("K122").Value = 3617.8287
$ws.Range("M122").Value = -1167.8287
$ws.Range("H126").Value = 3749.5833
$ws.Range("I126").Value = 2212
$ws.Range("J126").Value = 6824.75
$ws.Range("K126").Value = 6636
$ws.Range("L126").Value = 20474.25
$ws.Range("M126").Value = -4166
$ws.Range("N126").Value = -25414.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 18149.334
$ws.Range("I9").Value = 15000
$ws.Range("J9").Value = 19724
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 19724
$ws.Range("M9").Value = -14776
$ws.Range("N9").Value = -20172
$ws.Range("H40").Value = 6735.676
$ws.Range("I40").Value = 3621.913
$ws.Range("K40").Value = 3621.913
$ws.Range("M40").Value = -3485.913
$ws.Range("H61").Value = 3721.6667
$ws.Range("I61").Value = 3996.7778
$ws.Range("J61").Value = 2896.3333
$ws.Range("K61").Value = 3996.7778
$ws.Range("L61").Value = 2896.3333
$ws.Range("M61").Value = -3794.7778
$ws.Range("N61").Value = -3300.3333
$ws.Range("H113").Value = 3721.6667
$ws.Range("I113").Value = 3996.7778
$ws.Range("J113").Value = 2896.3333
$ws.Range("K113").Value = 3996.7778
$ws.Range("L113").Value = 2896.3333
$ws.Range("M113").Value = -1826.7778
$ws.Range("N113").Value = -7236.3333
$ws.Range("H122").Value = 3489.9048
$ws.Range("I122").Value = 3148.4856
$ws.Range("J122").Value = 5197
$ws.Range("K122").Value = 9445.4568
$ws.Range("L122").Value = 15591
$ws.Range("M122").Value = -6995.4568
$ws.Range("N122").Value = -20491
$ws.Range("H132").Value = 36730.047
$ws.Range("I132").Value = 45820.707
$ws.Range("J132").Value = 5821.8
$ws.Range("K132").Value = 137462.121
$ws.Range("L132").Value = 17465.4
$ws.Range("M132").Value = -134932.121
$ws.Range("N132").Value = -22525.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 57190.168
$ws.Range("I41").Value = 48547.332
$ws.Range("J41").Value = 65833
$ws.Range("K41").Value = 48547.332
$ws.Range("L41").Value = 65833
$ws.Range("M41").Value = -48157.332
$ws.Range("N41").Value = -66613
$ws.Range("H45").Value = 16935.625
$ws.Range("I45").Value = 15299.5
$ws.Range("J45").Value = 17481
$ws.Range("K45").Value = 15299.5
$ws.Range("L45").Value = 17481
$ws.Range("M45").Value = -14808.5
$ws.Range("N45").Value = -18463
$ws.Range("H122").Value = 7272.1113
$ws.Range("J122").Value = 7159.8
$ws.Range("L122").Value = 21479.4
$ws.Range("N122").Value = -26379.4
$ws.Range("H126").Value = 3995
$ws.Range("J126").Value = 3995
$ws.Range("L126").Value = 11985
$ws.Range("N126").Value = -16925
$ws.Range("H132").Value = 3132.7837
$ws.Range("I132").Value = 2872.3125
$ws.Range("J132").Value = 4799.8
$ws.Range("K132").Value = 8616.9375
$ws.Range("L132").Value = 14399.4
$ws.Range("M132").Value = -6086.9375
$ws.Range("N132").Value = -19459.4
$ws.Range("H136").Value = 5354.5654
$ws.Range("I136").Value = 3878.5293
$ws.Range("J136").Value = 9536.666999999999
$ws.Range("K136").Value = 11635.5879
$ws.Range("L136").Value = 28610.001
$ws.Range("M136").Value = -9085.5879
$ws.Range("N136").Value = -33710.001

Write-Host "Updated $([int]275) cells across 8 sheets"
